$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.027.85"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "2.274.40"
$ws.Range("E3").Value = "  +2.51%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.56"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.23"
$ws.Range("E6").Value = "  +5.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.588"
$ws.Range("E7").Value = "  +1.53%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.572"
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.52"
$ws.Range("E10").Value = "  +5.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0837"
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.86"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").Value = "2.621.67"
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.874"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.52"
$ws.Range("E16").Value = "  +3.66%  "
$ws.Range("D17").Value = "2.279.12"
$ws.Range("E17").Value = "  +2.81%  "
$ws.Range("D18").Value = "43.938.93"
$ws.Range("E18").Value = "  +2.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.33"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").Value = "0.0₃0996"
$ws.Range("E20").Value = "  +3.55%  "
$ws.Range("E21").Value = "  +3.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.12"
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("E25").Value = "  +4.28%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.26"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.53"
$ws.Range("E28").Value = "  +17.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.51"
$ws.Range("E30").Value = "  +4.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "162.21"
$ws.Range("E31").Value = "  +4.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.48"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0877"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.28"
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.53"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("E39").Value = "  +4.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.85"
$ws.Range("E40").Value = "  +5.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.70"
$ws.Range("E41").Value = "  +29.55%  "
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").Value = "1.779.21"
$ws.Range("E44").Value = "  -5.16%  "
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.74"
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "59.49"
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.73"
$ws.Range("E50").Value = "  -4.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.49"
$ws.Range("E51").Value = "  +3.84%  "
